# Applies the "Updated cryptos list" data refresh (Sat Sep 30 12:15:05 UTC 2023).
# Updates Price (D) / Volume(1h) (E) cells for existing coins, and swaps the
# ImmutableX / VeChain rows (37-38) to their new ranking + refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text looks like a plain number (e.g. "215.06") get a
# leading apostrophe so Excel stores them as text, exactly like the original sheet
# (prices such as "26.974.99" use "." as a thousands separator, so they are text,
# never real numbers).
$ws.Range("D2").Value = '26.974.99'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '1.678.57'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '''215.06'
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("E6").Value = '  -3.58%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.66%  '

$ws.Range("E9").Value = '  +0.22%  '

$ws.Range("D10").Value = '''20.35'
$ws.Range("E10").Value = '  +1.30%  '

$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").Value = '1.914.57'
$ws.Range("E12").Value = '  +0.51%  '

$ws.Range("D13").Value = '1.716.44'
$ws.Range("E13").Value = '  +2.79%  '

$ws.Range("E14").Value = '  +0.38%  '

$ws.Range("E15").Value = '  +1.68%  '

$ws.Range("D16").Value = '''65.73'
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").Value = '''8.21'
$ws.Range("E17").Value = '  +6.46%  '

$ws.Range("D18").Value = '27.008.79'
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").Value = '''235.77'
$ws.Range("E19").Value = '  +1.16%  '

$ws.Range("E20").Value = '  -0.23%  '

$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("D22").Value = '''4.44'
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("E23").Value = '  -0.34%  '

$ws.Range("E24").Value = '  -3.00%  '

$ws.Range("D25").Value = '''146.38'
$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("E26").Value = '  +1.02%  '

$ws.Range("D27").Value = '''16.05'
$ws.Range("E27").Value = '  +1.32%  '

$ws.Range("E28").Value = '  -3.96%  '

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").Value = '''0.0498'
$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("D33").Value = '1.482.13'
$ws.Range("E33").Value = '  +2.44%  '

$ws.Range("E34").Value = '  +1.45%  '

$ws.Range("D35").Value = '''1.69'
$ws.Range("E35").Value = '  +5.43%  '

$ws.Range("D36").Value = '''2.41'
$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.586'
$ws.Range("E37").Value = '  +3.55%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.0175'
$ws.Range("E38").Value = '  +3.67%  '

$ws.Range("D39").Value = '''0.905'
$ws.Range("E39").Value = '  +1.61%  '

$ws.Range("D40").Value = '''5.78'
$ws.Range("E40").Value = '  -4.56%  '

$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").Value = '''1.01'
$ws.Range("E42").Value = '  +0.26%  '

$ws.Range("E43").Value = '  +1.10%  '

$ws.Range("D44").Value = '''67.49'
$ws.Range("E44").Value = '  +2.98%  '

$ws.Range("D45").Value = '1.819.52'
$ws.Range("E45").Value = '  +0.24%  '

$ws.Range("D46").Value = '''0.783'
$ws.Range("E46").Value = '  +0.52%  '

$ws.Range("D47").Value = '''90.55'
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("E48").Value = '  +0.35%  '

$ws.Range("E49").Value = '  -0.32%  '

$ws.Range("E50").Value = '  +2.11%  '

$ws.Range("D51").Value = '''0.0508'
$ws.Range("E51").Value = '  +0.15%  '
